# Update column G ("K" strikeout count) values on Sheet1 to reflect the
# regenerated save data (switching from Strike# to K, recalculated s_vals).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 0
    3  = 4
    4  = 4
    5  = 3
    6  = 7
    7  = 4
    8  = 6
    9  = 7
    10 = 4
    11 = 3
    12 = 8
    13 = 3
    14 = 8
    15 = 8
    16 = 9
    17 = 8
    18 = 5
    19 = 5
    20 = 9
    21 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
